# Latest Update - Hate Crime
# Added Hate Crime as a measure: appends a new "Hate_Crime" worksheet
# (mirroring the layout of the other measure sheets) after "Theft_Person".

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# --- Add the new sheet after the last existing sheet (Theft_Person) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Hate_Crime"

# --- Header row ---
$ws.Range("A1").Value = "Month_Year"
$ws.Range("B1").Value = "count"

# --- Legend / chart metadata block (F2:G4), matching the other sheets ---
$ws.Range("F2").Value = "metal"
$ws.Range("G2").Value = "value"
$ws.Range("F3").Value = "Title"
$ws.Range("G3").Value = "Hate Crime Offences"
$ws.Range("F4").Value = "Subtitle"
$ws.Range("G4").Value = "Number of hate crime offences in London as recorded by the MPS"

# --- Copy the date-format style used by the other measure sheets onto column A ---
$styleSource = $wb.Worksheets.Item("Theft_Person")
$styleSource.Range("A2").Copy()
$ws.Range("A2:A71").PasteSpecial($xlPasteFormats)

# --- Monthly data: Month_Year (date serial) / count ---
$monthData = @(
    "42826|1610","42856|1705","42887|2149","42917|1988","42948|1757","42979|1538",
    "43009|1644","43040|1396","43070|1270","43101|1326","43132|1243","43160|1616",
    "43191|1561","43221|1749","43252|1820","43282|1932","43313|1496","43344|1506",
    "43374|1536","43405|1653","43435|1573","43466|1445","43497|1548","43525|1853",
    "43556|1774","43586|1806","43617|1874","43647|2077","43678|1908","43709|1784",
    "43739|1876","43770|1875","43800|1768","43831|1646","43862|1806","43891|1795",
    "43922|1472","43952|1873","43983|2760","44013|2710","44044|2616","44075|2286",
    "44105|1984","44136|1847","44166|1644","44197|1508","44228|1540","44256|2106",
    "44287|2131","44317|2573","44348|2537","44378|2810","44409|2219","44440|2188",
    "44470|2310","44501|2212","44531|2006","44562|1915","44593|1879","44621|2227",
    "45017|2002","45047|2305","45078|2425","45108|2330","45139|2242","45170|2305",
    "45200|2876","45231|2705","45261|2281","45292|2042"
)

$row = 2
foreach ($entry in $monthData) {
    $parts = $entry.Split("|")
    $ws.Cells.Item($row, 1).Value = [double]$parts[0]
    $ws.Cells.Item($row, 2).Value = [double]$parts[1]
    $row = $row + 1
}
